$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4194.8335
$ws.Range("I64").Value = 4037.6667
$ws.Range("J64").Value = 4666.3335
$ws.Range("K64").Value = 4037.6667
$ws.Range("L64").Value = 4666.3335
$ws.Range("M64").Value = -3789.6667
$ws.Range("N64").Value = -5162.3335
$ws.Range("H67").Value = 4194.8335
$ws.Range("I67").Value = 4037.6667
$ws.Range("J67").Value = 4666.3335
$ws.Range("K67").Value = 4037.6667
$ws.Range("L67").Value = 4666.3335
$ws.Range("M67").Value = -3179.6667
$ws.Range("N67").Value = -6382.3335
$ws.Range("H70").Value = 1035.8572
$ws.Range("J70").Value = 834
$ws.Range("L70").Value = 2502
$ws.Range("N70").Value = -3042
$ws.Range("H73").Value = 1035.8572
$ws.Range("J73").Value = 834
$ws.Range("L73").Value = 2502
$ws.Range("N73").Value = -4374
$ws.Range("H74").Value = 6253795
$ws.Range("I74").Value = 2634.3333
$ws.Range("J74").Value = 7356941
$ws.Range("K74").Value = 2634.3333
$ws.Range("L74").Value = 7356941
$ws.Range("M74").Value = -1698.3333
$ws.Range("N74").Value = -7358813
$ws.Range("H77").Value = 6253795
$ws.Range("I77").Value = 2634.3333
$ws.Range("J77").Value = 7356941
$ws.Range("K77").Value = 13171.6665
$ws.Range("L77").Value = 36784705
$ws.Range("M77").Value = -8491.666499999999
$ws.Range("N77").Value = -36794065
$ws.Range("H106").Value = 918.95
$ws.Range("I106").Value = 918.95
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 918.95
$ws.Range("L106").Value = 0
$ws.Range("M106").ClearContents() | Out-Null
$ws.Range("N106").Value = -287.95
$ws.Range("H129").Value = 770.81134
$ws.Range("I129").Value = 402.875
$ws.Range("J129").Value = 836.2222
$ws.Range("K129").Value = 1208.625
$ws.Range("L129").Value = 2508.6666
$ws.Range("M129").Value = 3791.375
$ws.Range("N129").Value = -12508.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3703.8572
$ws.Range("I32").Value = 3225.7917
$ws.Range("J32").Value = 6572.25
$ws.Range("K32").Value = 3225.7917
$ws.Range("L32").Value = 6572.25
$ws.Range("M32").Value = -2938.7917
$ws.Range("N32").Value = -7146.25
$ws.Range("H45").Value = 3269.111
$ws.Range("I45").Value = 2623.7
$ws.Range("J45").Value = 4075.875
$ws.Range("K45").Value = 2623.7
$ws.Range("L45").Value = 4075.875
$ws.Range("M45").Value = -2246.7
$ws.Range("N45").Value = -4829.875
$ws.Range("H63").Value = 2606248.8
$ws.Range("I63").Value = 2331.7778
$ws.Range("K63").Value = 2331.7778
$ws.Range("M63").Value = -1645.7778
$ws.Range("H66").Value = 2606248.8
$ws.Range("I66").Value = 2331.7778
$ws.Range("K66").Value = 11658.889
$ws.Range("M66").Value = -8226.888999999999
$ws.Range("H114").Value = 44800
$ws.Range("J114").Value = 44800
$ws.Range("L114").Value = 44800
$ws.Range("N114").Value = -53478
$ws.Range("H132").Value = 14636.641
$ws.Range("I132").Value = 1494.5807
$ws.Range("J132").Value = 65562.125
$ws.Range("K132").Value = 4483.742099999999
$ws.Range("L132").Value = 196686.375
$ws.Range("M132").Value = -1953.742099999999
$ws.Range("N132").Value = -201746.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 270.3
$ws.Range("I22").Value = 272.55554
$ws.Range("K22").Value = 272.55554
$ws.Range("M22").Value = -99.55554000000001
$ws.Range("H107").Value = 740
$ws.Range("I107").Value = 740
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 740
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents() | Out-Null
$ws.Range("N107").Value = 1180
$ws.Range("H134").Value = 2890.1614
$ws.Range("I134").Value = 2919.8333
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 8759.499899999999
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = -6224.499899999999
$ws.Range("N134").Value = -11070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3029.8057
$ws.Range("I31").Value = 2359.3
$ws.Range("J31").Value = 3867.9375
$ws.Range("K31").Value = 2359.3
$ws.Range("L31").Value = 3867.9375
$ws.Range("M31").Value = -2064.3
$ws.Range("N31").Value = -4457.9375
$ws.Range("H34").Value = 3029.8057
$ws.Range("I34").Value = 2359.3
$ws.Range("J34").Value = 3867.9375
$ws.Range("K34").Value = 2359.3
$ws.Range("L34").Value = 3867.9375
$ws.Range("M34").Value = -2157.3
$ws.Range("N34").Value = -4271.9375
$ws.Range("H122").Value = 2289.2222
$ws.Range("I122").Value = 2289.2222
$ws.Range("K122").Value = 6867.6666
$ws.Range("M122").Value = -4417.6666
$ws.Range("H134").Value = 1057.4615
$ws.Range("I134").Value = 767.3684
$ws.Range("K134").Value = 2302.1052
$ws.Range("M134").Value = 232.8948

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 5631.5
$ws.Range("J63").Value = 5631.5
$ws.Range("L63").Value = 16894.5
$ws.Range("N63").Value = -18392.5
$ws.Range("H66").Value = 5631.5
$ws.Range("J66").Value = 5631.5
$ws.Range("L66").Value = 50683.5
$ws.Range("N66").Value = -58171.5
$ws.Range("H74").Value = 7825
$ws.Range("I74").Value = 5800
$ws.Range("K74").Value = 17400
$ws.Range("M74").Value = -16339
$ws.Range("H77").Value = 7825
$ws.Range("I77").Value = 5800
$ws.Range("K77").Value = 52200
$ws.Range("M77").Value = -46896
$ws.Range("H80").Value = 2958.6
$ws.Range("J80").Value = 2948.25
$ws.Range("L80").Value = 8844.75
$ws.Range("N80").Value = -10716.75
$ws.Range("H83").Value = 2958.6
$ws.Range("J83").Value = 2948.25
$ws.Range("L83").Value = 26534.25
$ws.Range("N83").Value = -35894.25
$ws.Range("H113").Value = 14799.786
$ws.Range("J113").Value = 599.6667
$ws.Range("L113").Value = 1799.0001
$ws.Range("N113").Value = -6139.0001
$ws.Range("H121").Value = 13076.75
$ws.Range("J121").Value = 25674.75
$ws.Range("L121").Value = 77024.25
$ws.Range("N121").Value = -79644.25
$ws.Range("H131").Value = 108335.05
$ws.Range("J131").Value = 110694.055
$ws.Range("L131").Value = 332082.165
$ws.Range("N131").Value = -342162.165

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1968.8572
$ws.Range("I97").Value = 630.3333
$ws.Range("J97").Value = 10000
$ws.Range("K97").Value = 630.3333
$ws.Range("L97").Value = 10000
$ws.Range("M97").Value = -134.3333
$ws.Range("N97").Value = -10992
$ws.Range("H102").Value = 1615.8649
$ws.Range("I102").Value = 1674.4193
$ws.Range("J102").Value = 1313.3334
$ws.Range("K102").Value = 1674.4193
$ws.Range("L102").Value = 1313.3334
$ws.Range("M102").Value = -52.41930000000002
$ws.Range("N102").Value = -4557.3334
$ws.Range("H113").Value = 4280
$ws.Range("I113").Value = 2700
$ws.Range("K113").Value = 2700
$ws.Range("M113").Value = -530

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3083.5908
$ws.Range("I7").Value = 3133.6875
$ws.Range("J7").Value = 2950
$ws.Range("K7").Value = 3133.6875
$ws.Range("L7").Value = 2950
$ws.Range("M7").Value = -3021.6875
$ws.Range("N7").Value = -3174
$ws.Range("H46").Value = 3000
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents() | Out-Null
$ws.Range("H55").Value = 215.92857
$ws.Range("I55").Value = 193
$ws.Range("J55").Value = 222.18182
$ws.Range("K55").Value = 193
$ws.Range("L55").Value = 222.18182
$ws.Range("M55").Value = -20
$ws.Range("N55").Value = -568.18182
$ws.Range("H68").Value = 4178.8823
$ws.Range("I68").Value = 2138.2222
$ws.Range("J68").Value = 6474.625
$ws.Range("K68").Value = 2138.2222
$ws.Range("L68").Value = 6474.625
$ws.Range("M68").Value = -1389.2222
$ws.Range("N68").Value = -7972.625
$ws.Range("H71").Value = 4178.8823
$ws.Range("I71").Value = 2138.2222
$ws.Range("J71").Value = 6474.625
$ws.Range("K71").Value = 10691.111
$ws.Range("L71").Value = 32373.125
$ws.Range("M71").Value = -6947.111000000001
$ws.Range("N71").Value = -39861.125
$ws.Range("H93").Value = 2730.6
$ws.Range("I93").Value = 2867.3333
$ws.Range("J93").Value = 1500
$ws.Range("K93").Value = 2867.3333
$ws.Range("L93").Value = 1500
$ws.Range("M93").Value = -1619.3333
$ws.Range("N93").Value = -3996
$ws.Range("H126").Value = 3083.5908
$ws.Range("I126").Value = 3133.6875
$ws.Range("J126").Value = 2950
$ws.Range("K126").Value = 9401.0625
$ws.Range("L126").Value = 8850
$ws.Range("M126").Value = -6931.0625
$ws.Range("N126").Value = -13790
$ws.Range("H132").Value = 3111.889
$ws.Range("I132").Value = 2251
$ws.Range("J132").Value = 3800.6
$ws.Range("K132").Value = 6753
$ws.Range("L132").Value = 11401.8
$ws.Range("M132").Value = -4223
$ws.Range("N132").Value = -16461.8
$ws.Range("H136").Value = 1071.6875
$ws.Range("I136").Value = 1124.7858
$ws.Range("J136").Value = 700
$ws.Range("K136").Value = 3374.3574
$ws.Range("L136").Value = 2100
$ws.Range("M136").Value = -824.3574000000003
$ws.Range("N136").Value = -7200
